# Update the "Correspond Handoff Datetime" (column E) and
# "Correspond Handback DateTime" (column H) timestamps on row 4
# (the 931b8af4-... entry) for both the zh-cn and de-de sheets,
# reflecting a newly generated handback report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-19 08:36:20"
$wsZhCn.Range("H4").Value = "2016-03-19 08:36:40"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-19 08:36:23"
$wsDeDe.Range("H4").Value = "2016-03-19 08:36:46"
